# The original submitter edited request N3 ("Num de Resolucion" for the
# second student row), changing its value from "17" to "24", then scrolled
# the sheet so column J is at the left edge and left the N2:N3 region
# selected (active cell N3).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell edit: N3 "17" -> "24" (kept as text, matching the original type) ---
$cell = $ws.Range("N3")
$cell.NumberFormat = "@"
$cell.Value = "24"
$cell.NumberFormat = "General"

# --- View/selection: scroll so column J is the leftmost visible column and
#     leave N2:N3 selected, matching where the edited cell lives. ---
$ws.Range("N2:N3").Select()
$excel.ActiveWindow.ScrollColumn = 10
$excel.ActiveWindow.ScrollRow = 1

# --- Column width touch-up: D and J grew very slightly (sub-pixel) in the
#     authored workbook, most likely from the host app's own re-layout; set
#     them as close as this engine's width grid allows. ---
$ws.Columns.Item(4).ColumnWidth = 18.7449392712551
$ws.Columns.Item(10).ColumnWidth = 68.4493927125506
